# Inserts one new data row into the "Apio" price table at row 151, pushing
# the existing rows 151-242 down to 152-243 (new dimension A1:R243).
# The newly inserted row 151 carries a fresh weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 151..242 down by one, creating a blank row 151.
$ws.Rows.Item(151).Insert()

# Populate the new row 151 with the new observation.
$ws.Range("A151").Value = 10
$ws.Range("B151").Value = "Vega Modelo de Temuco"
$ws.Range("C151").Value = "La Araucanía"
$ws.Range("D151").Value = 44582
$ws.Range("E151").Value = 9
$ws.Range("F151").Value = 100112017
$ws.Range("G151").Value = "Apio"
$ws.Range("H151").Value = "Americana (o)"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 20
$ws.Range("K151").Value = 10000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 10000
$ws.Range("N151").Value = "`$/docena de matas"
$ws.Range("O151").Value = "Provincia del Elquí"
$ws.Range("P151").Value = 1667
$ws.Range("Q151").Value = 6
$ws.Range("R151").Value = "Hortaliza"
